# Update the "想去人数" (number of people interested) column (F) values
# for several events across the workbook's sheets, reflecting newer
# generated data (gh-pages output regenerated at commit 7921097).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 11460   # 广州·《FGO》FES2024·广州特别纪念展
$ws.Range("F5").Value = 859     # 广州·第七届萌物语动漫嘉年华
$ws.Range("F16").Value = 564    # 广州·运动番only4.0
$ws.Range("F18").Value = 1153   # 广州·樱漫动漫嘉年华8.0
$ws.Range("F19").Value = 220    # 广州·星火.AI动漫嘉年华3.0
$ws.Range("F28").Value = 513    # 广州·进击的巨人only
$ws.Range("F29").Value = 698    # 广州·代号鸢only2.0

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 78      # 广州·「十年之约」封茗囧菌2024个唱

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 96      # 广州·次元波板糖×线条小狗MALTESE 主题快闪店

# --- Sheet: 全部类型 (All types, aggregated) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 11460   # 广州·《FGO》FES2024·广州特别纪念展
$ws.Range("F6").Value = 859     # 广州·第七届萌物语动漫嘉年华
$ws.Range("F13").Value = 96     # 广州·次元波板糖×线条小狗MALTESE 主题快闪店
$ws.Range("F21").Value = 564    # 广州·运动番only4.0
$ws.Range("F23").Value = 1153   # 广州·樱漫动漫嘉年华8.0
$ws.Range("F24").Value = 220    # 广州·星火.AI动漫嘉年华3.0
$ws.Range("F31").Value = 78     # 广州·「十年之约」封茗囧菌2024个唱
$ws.Range("F38").Value = 513    # 广州·进击的巨人only
$ws.Range("F39").Value = 698    # 广州·代号鸢only2.0
